# MRCY Quarterly Financials - add two new quarterly columns (Dec-2018 and
# Sep-2018) in front of the existing quarter columns.
#
# This mirrors what happened in Excel: two new columns were inserted right
# before column D (shifting the old D:K data right to F:M), the inserted
# columns picked up the formatting of the columns immediately to their
# right, and then the new quarter's figures were typed into the now-empty
# D:E columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MRCY")

# 1) Insert two blank columns at D:E - this shifts the existing D:K data to F:M.
$ws.Range("D1:E1").EntireColumn.Insert()

# 2) The freshly inserted D:E columns don't carry the right number formats yet
#    (Excel seeds them from the column that used to be there). Copy the
#    formatting for the whole data block from F:G (the old D:E, now shifted)
#    onto the new D:E so dates/numbers keep displaying correctly.
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 3) Fill in the new quarter's values (Dec-2018 in D, Sep-2018 in E) for every
#    row that has data, row by row.
$newData = @(
  @(7, 43465, 43373),
  @(8, 159100, 144100),
  @(9, 88200, 82500),
  @(10, 70900, 61600),
  @(12, 16200, 14900),
  @(13, 0, 0),
  @(14, 100, 900),
  @(15, 6900, 7200),
  @(17, 139200, 130200),
  @(18, 19900, 13900),
  @(20, -800, -1000),
  @(21, 30800, 24400),
  @(22, 2200, 2300),
  @(23, 16900, 10600),
  @(24, 4500, 3100),
  @(25, 0, 0),
  @(26, 12400, 7500),
  @(27, 12400, 7500),
  @(28, 0, 0),
  @(29, "NA", "NA"),
  @(30, 0, 0),
  @(31, 0, 0),
  @(32, 800, 1000),
  @(33, 12400, 7500),
  @(34, 0, 0),
  @(35, 12400, 7500),
  @(38, 43465, 43373),
  @(41, 93900, 72900),
  @(42, 0, 0),
  @(43, 168300, 153900),
  @(44, 126400, 121200),
  @(45, 10700, 16300),
  @(46, 399300, 364200),
  @(47, 0, 0),
  @(48, 53100, 50800),
  @(49, 696300, 704200),
  @(50, 0, 0),
  @(51, 0, 0),
  @(52, 7900, 7800),
  @(53, 0, 0),
  @(54, 1156600, 1127000),
  @(57, 30800, 25700),
  @(58, "NA", "NA"),
  @(59, 62400, 57000),
  @(60, 93200, 82800),
  @(61, 240000, 240000),
  @(62, 27400, 27600),
  @(63, 0, 0),
  @(64, 0, 0),
  @(65, 0, 0),
  @(66, 360600, 350400),
  @(68, 0, 0),
  @(69, 0, 0),
  @(70, 0, 0),
  @(71, 0, 0),
  @(72, 199800, 187400),
  @(73, 0, 0),
  @(74, 0, 0),
  @(75, 0, 0),
  @(76, 796100, 776600),
  @(77, 0, 0),
  @(80, 43465, 43373),
  @(81, 12400, 7500),
  @(83, 11700, 11500),
  @(84, 0, 0),
  @(85, 0, 0),
  @(86, 0, 0),
  @(87, 0, 0),
  @(88, 0, 0),
  @(89, 25300, 20000),
  @(91, -7100, -3700),
  @(92, 0, 0),
  @(93, 0, 0),
  @(94, -5800, -50000),
  @(96, 0, 0),
  @(97, 0, 0),
  @(98, 0, 0),
  @(99, 0, 0),
  @(100, 1600, 36300),
  @(101, 0, 0),
  @(102, 21000, 6300)
)

foreach ($item in $newData) {
  $r = $item[0]
  $ws.Cells.Item($r, 4).Value = $item[1]
  $ws.Cells.Item($r, 5).Value = $item[2]
}

Write-Output "done"
